$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 21.28240740740739
$ws.Range("N2").Value = 2.017497406510892
$ws.Range("O2").Value = 2.228623569098047

$ws.Range("I6").Value = 19.65277777777778
$ws.Range("N6").Value = 1.983015294974508
$ws.Range("O6").Value = 2.18606997558991

$ws.Range("I7").Value = 13.75752314814816
$ws.Range("N7").Value = 1.867546171126113
$ws.Range("O7").Value = 2.044826120875009

$ws.Range("I10").Value = 19.65277777777778
$ws.Range("N10").Value = 1.983015294974508
$ws.Range("O10").Value = 2.18606997558991
